$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.630.28"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "1.893.26"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "'310.15"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "'0.5238"
$ws.Range("E7").Value = "  +4.54%  "

$ws.Range("D8").Value = "'0.3806"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'0.07231"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("D10").Value = "'21.05"
$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("D11").Value = "'0.9000"

$ws.Range("D12").Value = "1.883.46"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").Value = "'0.07623"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "'5.423"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "'0.000008666"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").Value = "'1.0000"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").Value = "27.678.02"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").Value = "'5.150"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "2.119.35"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "'10.79"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'6.590"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").Value = "'153.09"
$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("D26").Value = "'1.855"

$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").Value = "'2.186"
$ws.Range("E28").Value = "  -1.47%  "

$ws.Range("D29").Value = "'113.79"
$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("D30").Value = "'4.823"
$ws.Range("E30").Value = "  -1.98%  "

$ws.Range("D31").Value = "'4.786"
$ws.Range("E31").Value = "  +3.30%  "

$ws.Range("D32").Value = "'0.09107"
$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("D33").Value = "'0.05266"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "'3.116"
$ws.Range("E34").Value = "  -3.20%  "

$ws.Range("D35").Value = "'1.216"
$ws.Range("E35").Value = "  -1.62%  "

$ws.Range("D36").Value = "'0.7681"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").Value = "'0.02078"
$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("D38").Value = "'2.542"
$ws.Range("E38").Value = "  -0.41%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.091"
$ws.Range("E40").Value = "  -0.61%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5563"
$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("D42").Value = "'6.719"
$ws.Range("E42").Value = "  -3.48%  "

$ws.Range("D43").Value = "'116.39"
$ws.Range("E43").Value = "  +4.66%  "

$ws.Range("D44").Value = "'8.658"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("D45").Value = "'0.1505"

$ws.Range("D46").Value = "'0.4778"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("D47").Value = "'10.43"
$ws.Range("E47").Value = "  -1.81%  "

$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").Value = "'1.590"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("D50").Value = "'66.14"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("D51").Value = "'37.06"
$ws.Range("E51").Value = "  +0.24%  "
